$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2000
$ws.Range("B2").Value = 'PROJECTS VENDOR TEAM'
$ws.Range("C2").Value = 'If you are a vendor or outsourcer do send us your project and we outsource your projects across India. we do have a long list of satisfied customers data and across India we do have a tie-up with approx 500 vendors so if you have any new process send us the complete details along with SLA and we will outsource for you with no consultancy charges.'
$ws.Range("E2").Value = 'Projects Vendor is a Professional Services Firm'
$ws.Range("F2").Value = ' offshore talent, and enabling technologies, to optimize maximum utilization'
$ws.Range("G2").Value = 'https://www.linkedin.com/in/projectsvendor/?originalSubdomain=in'

$ws.Range("A3").Value = 2001
$ws.Range("B3").Value = 'JioMart - A New Giant in Digital Grocery Ecosystem'
$ws.Range("C3").Value = ': Grocery was 70%of Indian retail market with 90% driven by neighborhood shops. Reliance Industries had launched JioMart,an online grocery store in over 200 towns across the country. JioMart was expected to give a tough competition to prominent grocery delivery services like Amazon, Flipkart and Big Basket in the country. The national roll-out provided customers with a useful alternative to existing players (BigBasket, Grofers) as online grocery orders had spiked, reducing the store visits. Reliance Retail operated neighbourhood stores, supermarkets, hypermarkets, wholesale and online stores. The company had integrated its registered customer database from Reliance Fresh and other retail businesses into JioMart. Reliance’s presence in consumer electronics, fashion, lifestyle and grocery segments wasalready strong. Its retail business had crossedRs 45,000 crore in revenue in the December quarter of 2019'
$ws.Range("E3").Value = ': JioMart'
$ws.Range("F3").Value = ' offshore talent, and enabling technologies, to optimize maximum utilization'
$ws.Range("G3").Value = 'https://www.ijsr.net/archive/v9i10/SR201017180911.pdf'

$ws.Range("A4").Value = 2002
$ws.Range("B4").Value = 'Green Energy Project at Reliance AGM'
$ws.Range("C4").Value = 'The company had announced that it has a target to reach a net-zero carbon company by 2035, where 60% of the revenue that is earned by the company is through hydrocarbon-fueled energy operations, as per the reports of the recent financial year. While compared to other players in the industry 2035 is a relatively shorter time frame where the major players such as BPCL, Royal Dutch shell, and many others commit to reaching the target by 2050.'
$ws.Range("E4").Value = ':solar manufacturing unit'
$ws.Range("F4").Value = 'Giga Complex, which is expected to be one of the largest renewable energy integrated manufacturing facilities in the entire world.'
$ws.Range("G4").Value = 'https://startuptalky.com/reliance-agm-2021-highlights/'

$ws.Range("A5").Value = 2003
$ws.Range("B5").Value = 'Jio and Microsoft partnership by Reliance AGM'
$ws.Range("C5").Value = 'The company has also partnered with Microsoft in order to launch a 100 MW capacity of Jio-Azure cloud data centers in two cities of the country – Nagpur, and Jamnagar. The company is currently onboarding a group of customers for testing the product and service. Over the coming quarters, the company has plans to expand the data center capacity and offering the services to a growing number of SMEs and various other startups.'
$ws.Range("E5").Value = ': Reliance AGM'
$ws.Range("F5").Value = 'Green energy Giga factory that will be set up by the Reliance Industries'
$ws.Range("G5").Value = 'https://startuptalky.com/reliance-agm-2021-highlights/'
